$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 24 (event_id 14687083, synced 2025-09-20) ---
$ws.Range("A24").Value = "'14687083"
$ws.Range("B24").Value = "'2025-09-20"
$ws.Range("C24").Value = "Corentin Moutet"
$ws.Range("D24").Value = "Arthur Cazaux"
$ws.Range("E24").Value = "Gana Corentin Moutet"
$ws.Range("F24").Value = 2
$ws.Range("G24").Font.Bold = $false
$ws.Range("H24").Font.Bold = $false

# --- New row 25 (event_id 14721398, synced 2025-09-20) ---
$ws.Range("A25").Value = "'14721398"
$ws.Range("B25").Value = "'2025-09-20"
$ws.Range("C25").Value = "Petr Bar Biryukov"
$ws.Range("D25").Value = "Akira Santillan"
$ws.Range("E25").Value = "Gana Akira Santillan"
$ws.Range("F25").Value = 1.73
$ws.Range("G25").Font.Bold = $false
$ws.Range("H25").Font.Bold = $false
